# UserStory 7 tasks erstellt
#
# Adds four new task rows for UserStory 7 ("beim buchen die Preise
# einsehen und eine email als Bestätigung bekommen") to the Tasks sheet,
# and updates the remembered cell selections on the Priorisierung and
# Tasks sheets.

$wb = $excel.ActiveWorkbook

# --- Priorisierung sheet: just move the remembered selection ---
$prioSheet = $wb.Worksheets.Item("Priorisierung")
$null = $prioSheet.Range("B13").Select()

# --- Tasks sheet: add the four new task rows under UserStory 7 ---
$tasksSheet = $wb.Worksheets.Item("Tasks")

$tasksSheet.Range("B50").Value = "GetPrice in IBookingRepository und BookingRepositoryDB erstellen"
$tasksSheet.Range("B52").Value = "Price in Booking View einbauen"
$tasksSheet.Range("B53").Value = "Email Bestätigung bei BookingConfirmation einbauen"
$tasksSheet.Range("B51").Value = "Price in BookingController einbauen"

# Leave the Tasks sheet active with the new last entry selected.
$null = $tasksSheet.Range("B51").Select()
